# Applies the "Traps for young players" edit:
#  - Adds spell/grammar proofing marks (w:proofErr) around certain words,
#    splitting the affected runs accordingly.
#  - Appends several new bullet points (ListParagraph / numId 2) after the
#    "Make sure you use the correct target (hat vs racer)" bullet, including
#    a final, empty bullet.
#
# Because Word's COM object model has no direct "insert a proofErr mark"
# property, we rebuild the affected paragraphs (and add the new ones) by
# replacing/inserting their underlying OOXML via Range.InsertXML - this is
# still standard Range-based editing, just expressed as WordprocessingML
# fragments instead of plain text.

$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

function Set-ParagraphXml($paragraph, [string]$pAttrs, [string]$innerXml) {
    $range = $paragraph.Range
    $xml = '<w:p xmlns:w="' + $wNs + '" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"' + $pAttrs + '>' + $innerXml + '</w:p>'
    $range.InsertXML($xml)
}

function Assert-ParagraphStartsWith($paragraph, [string]$expectedPrefix) {
    if (-not ($paragraph.Range.Text.StartsWith($expectedPrefix))) {
        throw "Unexpected paragraph content. Expected prefix '$expectedPrefix' but found '$($paragraph.Range.Text)'"
    }
}

# --- 1. "OpenOCD not recognised? Launch VScode through ..." bullet -------
$p1 = $d.Paragraphs(2)
Assert-ParagraphStartsWith $p1 "OpenOCD not recognised"
$attrs1 = ' w14:paraId="57C73CC6" w14:textId="69489BFC" w:rsidR="00CD1526" w:rsidRDefault="00CD1526" w:rsidP="00CD1526"'
$inner1 = @'
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>OpenOCD</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> not recognised? Launch </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>VScode</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> through </w:t></w:r><w:r w:rsidRPr="00CD1526"><w:t>vscode.bat</w:t></w:r><w:r><w:t xml:space="preserve"> in wacky racers.</w:t></w:r>
'@
Set-ParagraphXml $p1 $attrs1 $inner1

# --- 2. "Did you run openOCD?" bullet -------------------------------------
$p2 = $d.Paragraphs(3)
Assert-ParagraphStartsWith $p2 "Did you run openOCD"
$attrs2 = ' w14:paraId="2194F11C" w14:textId="1B73182D" w:rsidR="0012377D" w:rsidRDefault="0012377D" w:rsidP="00CD1526"'
$inner2 = @'
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Did you run </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>openOCD</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>?</w:t></w:r>
'@
Set-ParagraphXml $p2 $attrs2 $inner2

# --- 3. "Board not working like it should? Check pins ..." bullet --------
$p3 = $d.Paragraphs(6)
Assert-ParagraphStartsWith $p3 "Board not working like it should"
$attrs3 = ' w14:paraId="2F40B25B" w14:textId="4803143D" w:rsidR="00CD1526" w:rsidRDefault="00CD1526" w:rsidP="00CD1526"'
$inner3 = @'
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Board not working like it </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>should?</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Check pins in the configuration file </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>target.h</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> are correct.</w:t></w:r>
'@
Set-ParagraphXml $p3 $attrs3 $inner3

# --- 4. New bullets appended after "Make sure you use the correct target
#        (hat vs racer)" -------------------------------------------------
$lastBullet = $d.Paragraphs(7)
Assert-ParagraphStartsWith $lastBullet "Make sure you use the correct target"
$newParasXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Make sure the correct silk screen label is used. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Eg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Level shifters silk screens underneath on hat.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>ST-Link does not provide power to anything other than the chip</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">USB does not provide LED Tape power </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Unconfigured LED\u2019s will glow dim</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr></w:p>
'@
# Swap in the real curly apostrophe (kept as an escape above to avoid any
# source-encoding ambiguity in this file).
$newParasXml = $newParasXml.Replace("\u2019", [char]0x2019)

$endOfLastBullet = $lastBullet.Range.End
$insertionRange = $d.Range($endOfLastBullet, $endOfLastBullet)
$insertionRange.InsertXML($newParasXml)
